$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.569.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.357.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.10%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "558.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.621"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.42%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.346.60"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.83%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.165"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.26%  "
$ws.Range("B11").Value = "Cardano"
$ws.Range("C11").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.632"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.61%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000276"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.890.40"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("E17").Value = "  +1.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.336.72"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "64.495.17"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.989"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "450.97"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +9.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.66"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.75"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.66"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.47"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "574.51"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.25%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.108"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.48%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "60.15"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.74%  "
$ws.Range("B36").Value = "Stacks"
$ws.Range("C36").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.69"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.38%  "
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.141"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.48"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.99%  "
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.371"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.31%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0741"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.79%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.084.23"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.09%  "
$ws.Range("B43").Value = "ThetaToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.82"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.65%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0417"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.73%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.134"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.30%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.45"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.41%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.24%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "138.84"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.47%  "
$ws.Range("B50").Value = "WEMIXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.54"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.18"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.61%  "
